$d = $word.ActiveDocument

# Build the WordprocessingML fragment for the new paragraph exactly as produced
# by the parser when it hits an unterminated "m:commentblock" user field: the
# field text is split across several runs ("{", "m", ":comment", "block",
# " some important comment", "}") followed by a bold red error-message run.
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$errMsg = "    &lt;---Invalid comment statement: Unexpected tag EOF missing [ENDCOMMENTBLOCK] while parsing m:commentblock some important comment"

$xml = '<w:p ' + $w + '>' `
  + '<w:r><w:t>{</w:t></w:r>' `
  + '<w:r><w:t>m</w:t></w:r>' `
  + '<w:r><w:t>:comment</w:t></w:r>' `
  + '<w:r><w:t>block</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> some important comment</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">}</w:t></w:r>' `
  + '<w:r><w:rPr><w:b w:val="on"/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">' + $errMsg + '</w:t></w:r>' `
  + '</w:p>'

$end = $d.Content.End
$r = $d.Range($end, $end)
$r.InsertXML($xml) | Out-Null
